$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Mac-Address rows (regcntr_id, machine_id, device_id) to append
# starting at row 147, matching the existing table layout:
# A=regcntr_id, B=machine_id, C=device_id, D=lang_code, E=is_active, F=cr_by, G=cr_dtimes
$newRows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update the selection to match the post-edit state
$ws.Range("H149").Select()
